$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that are stored as text (strings) in the workbook,
# even though they look numeric (e.g. "253.28"). Force text format on each cell
# individually before writing the value so Excel keeps it as a string instead of
# converting it to a numeric cell.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '253.28'
$ws.Range("D3").Value = '21.94'
$ws.Range("D4").Value = '5.567'
$ws.Range("D5").Value = '0.05704'
$ws.Range("D6").Value = '6.472'
$ws.Range("D7").Value = '0.8085'
$ws.Range("D8").Value = '1.044'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '0.1428'
$ws.Range("E9").Value = '8WazirXWRX'
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Value = '0.07326'
$ws.Range("E10").Value = '9MandalaExchangeTokenMDX'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '0.03147'
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B12").Value = 'ProBitToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D12").Value = '0.1266'
$ws.Range("E12").Value = '11ProBitTokenPROB'
$ws.Range("D13").Value = '0.02938'
$ws.Range("D14").Value = '0.09283'
$ws.Range("D15").Value = '0.001673'
$ws.Range("D16").Value = '3.219'
$ws.Range("D17").Value = '0.04764'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = '0.0005818'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").Value = '0.006455'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").Value = '0.005076'
$ws.Range("E20").Value = '19HotbitTokenHTBBestin24h'
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").Value = '0.001052'
$ws.Range("E21").Value = '20BitKanKAN'
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").Value = '0.0001500'
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").Value = '3.990'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("B24").Value = 'GateToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D24").Value = '3.384'
$ws.Range("E24").Value = '23GateTokenGT'
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").Value = '2.112'
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").Value = '0.3321'
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'
$ws.Range("D40").Value = '0.04144'
$ws.Range("D41").Value = '0.006919'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = '0.003499'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D43").Value = '0.1048'
$ws.Range("E43").Value = '42BKEXTokenBKK'
$ws.Range("D44").Value = '0.009570'
$ws.Range("D45").Value = '0.00005650'
$ws.Range("D47").Value = '0.7850'
$ws.Range("D48").Value = '0.01692'
